$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.000003350246346883397
$ws.Range("E2").Value = 0.000003350246346883397

# Row 3
$ws.Range("D3").Value = 0.9999993737879103
$ws.Range("E3").Value = 0.9999993737879103

# Row 4
$ws.Range("D4").Value = 0.001940857690868451
$ws.Range("E4").Value = 0.001940857690868451

# Row 5
$ws.Range("D5").Value = 0.000703243848620176
$ws.Range("E5").Value = 0.000703243848620176

# Row 6
$ws.Range("D6").Value = 0.2085961200085212
$ws.Range("E6").Value = 0.2085961200085212

# Row 7
$ws.Range("C7").Value = $false
$ws.Range("D7").Value = 0.0000002748387609522148
$ws.Range("E7").Value = 0.999999725161239

# Row 8
$ws.Range("D8").Value = 0.9666032262413636
$ws.Range("E8").Value = 0.03339677375863637

# Row 9
$ws.Range("D9").Value = 0.9909612988229372
$ws.Range("E9").Value = 0.009038701177062847

# Row 10
$ws.Range("D10").Value = 0.9999999999878302
$ws.Range("E10").Value = 0.00000000001216982070673112

# Row 11
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 2.967030048370361
$ws.Range("G11").Value = 0.8

# Row 12
$ws.Range("D12").Value = 0.00000007230541161122107
$ws.Range("E12").Value = 0.00000007230541161122107

# Row 13
$ws.Range("D13").Value = 0.00362969747672327
$ws.Range("E13").Value = 0.00362969747672327

# Row 14
$ws.Range("D14").Value = 0.0001793884047319373
$ws.Range("E14").Value = 0.0001793884047319373

# Row 15
$ws.Range("D15").Value = 0.00008059239550830151
$ws.Range("E15").Value = 0.00008059239550830151

# Row 16
$ws.Range("D16").Value = 0.07418202768473928
$ws.Range("E16").Value = 0.07418202768473928

# Row 17
$ws.Range("C17").Value = $false
$ws.Range("D17").Value = 0.0000000003931660192611175
$ws.Range("E17").Value = 0.9999999996068339

# Row 18
$ws.Range("D18").Value = 0.9951765220807354
$ws.Range("E18").Value = 0.004823477919264629

# Row 19
$ws.Range("D19").Value = 0.9995181554110644
$ws.Range("E19").Value = 0.0004818445889356315

# Row 20
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0

# Row 21
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 2.174308061599731
$ws.Range("G21").Value = 0.9
